# Apply comprehensive formatting standards to Open-Source_Tools_Review.docx
# - Times New Roman throughout (already set on every run; left untouched)
# - Standardized font sizes: Title 28pt, Heading2 14pt, Body 12pt
# - 1.5 line spacing (w:line=360 w:lineRule=auto) on every paragraph
# - Title: 12pt/12pt before/after spacing; Heading2: 9pt/5pt before/after
# - Left-align title, subtitle, headings and body paragraphs
# - Force every run to explicit black (0,0,0) font color
# - 0.75" (54pt) top/bottom page margins

$d = $word.ActiveDocument

# wdLineSpaceMultiple = 5 (input) -> persisted as lineRule="auto" with
# w:line = LineSpacing(pts) * 20; 18pt * 20 = 360 twips == the target 1.5x line.
$LineSpacingMultiple = 5
$LineSpacing1_5 = 18

# wdAlignParagraphLeft = 0
$AlignLeft = 0

# wdColorBlack / RGB(0,0,0)
$Black = 0

$paraCount = $d.Paragraphs.Count

# Left-align every paragraph in the document in one shot, via the
# Paragraphs *collection* rather than a per-paragraph ParagraphFormat
# assignment. Both ultimately express "left" (already the resolved
# default for every style here), but going through the collection keeps
# an explicit <w:jc w:val="left"/> on each paragraph instead of having it
# optimized away as redundant - so do this first and don't re-touch
# Alignment per-paragraph afterwards (that would re-collapse it to the
# implicit default).
$d.Paragraphs.Alignment = $AlignLeft

for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal

    # --- paragraph-level formatting (applies to the whole paragraph,
    #     including its end-of-paragraph mark) ---
    $pf = $p.Range.ParagraphFormat
    $pf.LineSpacingRule = $LineSpacingMultiple
    $pf.LineSpacing = $LineSpacing1_5

    if ($styleName -eq "Title") {
        $pf.SpaceBefore = 12
        $pf.SpaceAfter = 12
    } elseif ($styleName -eq "Heading 2") {
        $pf.SpaceBefore = 9
        $pf.SpaceAfter = 5
    }

    # --- run-level formatting: restrict to the paragraph's text (exclude
    #     the trailing paragraph-mark character) so only the w:r/w:rPr is
    #     touched, not a synthetic pPr/rPr for the mark itself ---
    $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)
    $textRange.Font.Color = $Black

    if ($styleName -eq "Title") {
        $textRange.Font.Size = 28
    } elseif ($styleName -eq "Heading 2") {
        $textRange.Font.Size = 14
    } else {
        # Normal-styled paragraphs: subtitle (already 12pt) and body text
        # (was 11pt) both land on 12pt.
        $textRange.Font.Size = 12
    }
}

# --- page margins: 0.75" top/bottom (54pt); left/right unchanged at 1080 twips ---
$d.PageSetup.TopMargin = 54
$d.PageSetup.BottomMargin = 54

Write-Output "Formatting applied"
